$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 59 (old rows 59-67 shift down to 61-69).
$ws.Rows("59:60").Insert()

# New row 59: Jengibre "Primera" entry for the week of 2021-11-22 (serial 44522)
$ws.Cells.Item(59, 1).Value = 9
$ws.Cells.Item(59, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(59, 3).Value = "Metropolitana"
$ws.Cells.Item(59, 4).Value = 44522
$ws.Cells.Item(59, 5).Value = 13
$ws.Cells.Item(59, 6).Value = 100114007
$ws.Cells.Item(59, 7).Value = "Jengibre"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 790
$ws.Cells.Item(59, 11).Value = 16000
$ws.Cells.Item(59, 12).Value = 18000
$ws.Cells.Item(59, 13).Value = 16987
$ws.Cells.Item(59, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(59, 15).Value = "Perú"
$ws.Cells.Item(59, 16).Value = 1307
$ws.Cells.Item(59, 17).Value = 13
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# New row 60: Jengibre "Segunda" entry for the week of 2021-11-22 (serial 44522)
$ws.Cells.Item(60, 1).Value = 9
$ws.Cells.Item(60, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60, 3).Value = "Metropolitana"
$ws.Cells.Item(60, 4).Value = 44522
$ws.Cells.Item(60, 5).Value = 13
$ws.Cells.Item(60, 6).Value = 100114007
$ws.Cells.Item(60, 7).Value = "Jengibre"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Segunda"
$ws.Cells.Item(60, 10).Value = 360
$ws.Cells.Item(60, 11).Value = 15000
$ws.Cells.Item(60, 12).Value = 15000
$ws.Cells.Item(60, 13).Value = 15000
$ws.Cells.Item(60, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(60, 15).Value = "Perú"
$ws.Cells.Item(60, 16).Value = 1154
$ws.Cells.Item(60, 17).Value = 13
$ws.Cells.Item(60, 18).Value = "Hortaliza"
